# Gates Demo Final
# Applies the authored changes to scan_childvacc_825_pg3.xlsx:
#  - Flip "hideInContents" (col F) to TRUE for the 8 "clause" rows that were missing it
#  - Widen the "values_list" column (col E) on the survey sheet
#  - Move the selection/active-tab from survey!F92 to choices!C8 (choices becomes active sheet)
#  - Convert the "choices" sheet's data_value / display.text columns (B:C, rows 2-7)
#    from the numeric literal 1 to the text value "yes"

$wb = $excel.ActiveWorkbook

$survey  = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# --- survey sheet: mark the begin-screen "clause" rows as hideInContents = TRUE ---
$hiddenRows = @(60, 64, 68, 72, 76, 80, 84, 88)
foreach ($r in $hiddenRows) {
    $survey.Range("F$r").Value = $true
}

# --- survey sheet: widen column E (values_list) ---
$survey.Columns.Item(5).ColumnWidth = 50

# --- choices sheet: change the data_value / display.text columns to "yes" ---
for ($r = 2; $r -le 7; $r++) {
    $choices.Range("B$r").Value = "yes"
    $choices.Range("C$r").Value = "yes"
}

# --- window / selection state ---
# Leave the survey sheet scrolled/selected first...
$survey.Activate() | Out-Null
$survey.Range("F40").Select() | Out-Null

# ...then make "choices" the active sheet/selection, matching the committed state.
$choices.Activate() | Out-Null
$choices.Range("C8").Select() | Out-Null

$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 0
$win.Width = 25600
$win.Height = 15240
